$d = $word.ActiveDocument

# 1. Capitalize "euclidean" -> "Euclidean" in the DI discussion paragraph.
$d.Content.Find.Execute("small network and euclidean distances", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "small network and Euclidean distances", 2) | Out-Null

# 2. Split the "We did so by assuming..." sentence into two bullet points,
#    breaking right after "network" and giving it its own closing period.
$d.Content.Find.Execute("on the network surrounding Lene-Voigt-Park.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "on the network.^psurrounding Lene-Voigt-Park.", 2) | Out-Null

# 3. Mark the "Tree / graphs like in Wolff, Scheuer et al. 2020" bullet as German text.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Tree / graphs like in Wolff, Scheuer et al. 2020") {
        $p.Range.LanguageID = "de-DE"
        break
    }
}

Write-Output "done"
